$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for two sightings (row 2 and row 3) had their Id ("A"),
# Antal ("I"), Ost ("Q") and Nord ("R") values swapped between each other.

# Row 2 becomes what used to be row 3's values
$ws.Range("A2").Value = 112098047
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "50"
$ws.Range("Q2").Value = 528942.9447194069
$ws.Range("R2").Value = 6229759.000311463

# Row 3 becomes what used to be row 2's values
$ws.Range("A3").Value = 112098022
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "10"
$ws.Range("Q3").Value = 528907.8083200558
$ws.Range("R3").Value = 6229763.767073607
